# LV_Activities - 26 June 2024
#
# Update the "James Craven" / "CapProviderTestCompany" placeholder values to
# "Amanda Donovan" / "ActivityCompany", then leave the selection/active-sheet
# state the way the author left it when they saved (Company sheet active,
# with a new selection on each sheet).

$wb = $excel.ActiveWorkbook

$users   = $wb.Worksheets.Item("Users")
$company = $wb.Worksheets.Item("Company")

# Update the data values.
$users.Range("A2").Value   = "Amanda Donovan"
$company.Range("A2").Value = "ActivityCompany"

# Move the selection on the (previously active) Users sheet.
$users.Range("C8").Select()

# Make Company the active sheet and move its selection too.
$company.Activate()
$company.Range("C11").Select()
